$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh - GitHub Actions scheduled update.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). D/E values are stored as
# text (matching the source feed formatting), so price cells are written
# with NumberFormat forced to Text to avoid Excel auto-converting numeric-
# looking strings ("46.71" etc.) into floating point numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '46.717.69'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.11%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.256.06'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.47%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '297.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.56%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.34%  '

# Row 7
$ws.Range("E7").Value = '  -1.14%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
$ws.Range("E9").Value = '  -6.27%  '

# Row 10
$ws.Range("E10").Value = '  -3.42%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0779'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.43%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.94'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.11%  '

# Row 13
$ws.Range("E13").Value = '  -1.81%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.594.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.69%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.254.23'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.79%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '46.622.75'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.08%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.42'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.13%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.786'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.08%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0963'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.70%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -10.32%  '

# Row 21
$ws.Range("E21").Value = '  -7.10%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.18'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.98%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '243.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.04%  '

# Row 24
$ws.Range("E24").Value = '  -7.20%  '

# Row 25
$ws.Range("E25").Value = '  -0.02%  '

# Row 26
$ws.Range("E26").Value = '  -7.42%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '40.48'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.16%  '

# Row 28
$ws.Range("E28").Value = '  -4.18%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.45'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.10%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.85'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.29%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.80'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.22%  '

# Row 32
$ws.Range("E32").Value = '  +3.41%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '143.36'
$ws.Range("D33").Style = "Normal"

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.27'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.39%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0760'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.95%  '

# Row 36
$ws.Range("E36").Value = '  -0.40%  '

# Row 37
$ws.Range("E37").Value = '  -3.18%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '15.23'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.76%  '

# Row 39
$ws.Range("E39").Value = '  -9.53%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.77'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.24%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0292'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.94%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.06'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -9.53%  '

# Row 43
$ws.Range("E43").Value = '  -0.03%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '92.31'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +13.67%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.778.33'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.40%  '

# Row 46
$ws.Range("E46").Value = '  -6.69%  '

# Row 47
$ws.Range("B47").Value = 'ordi'
$ws.Range("C47").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '69.66'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.11%  '

# Row 48
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.182'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.41%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.73'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.09%  '

# Row 50
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '93.11'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.56%  '

# Row 51
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.471.75'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.83%  '
